$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# "Shrink text on overflow" -> <a:bodyPr><a:normAutofit/></a:bodyPr>
$tf.AutoSize = 2

# Locate the two paragraphs we need to touch by their (pre-edit) text so the
# character offsets used by Characters(start, length) are derived rather
# than hard-coded.
$fullText = $tr.Text

# --- Paragraph 3 (lvl=2) ------------------------------------------------
# Old: "Integrated capabilities for applications"   (a single run)
# New: 5 runs -> "Integrated end-to-end solutio" / "n" /
#      " for applications requires transparent federation of  capabilities " /
#      "and the interoperability " / "of services"
#
# Edit this paragraph *before* paragraph 1 below: it sits later in the
# text, so shrinking/growing paragraph 1 afterwards can't invalidate the
# offsets already used here.
$oldP3 = "Integrated capabilities for applications"
$p3Start = $fullText.IndexOf($oldP3) + 1   # Characters() is 1-based
$run = $tr.Characters($p3Start, $oldP3.Length)

$newP3 = "Integrated end-to-end solution for applications requires transparent federation of  capabilities and the interoperability of services"
$run.Text = $newP3

$a = "Integrated end-to-end solutio"
$b = "n"
$c = " for applications requires transparent federation of  capabilities "
$d = "and the interoperability "
$e = "of services"

$pos = $p3Start
$tr.Characters($pos, $a.Length).Text = $a
$pos += $a.Length
$tr.Characters($pos, $b.Length).Text = $b
$pos += $b.Length
$tr.Characters($pos, $c.Length).Text = $c
$pos += $c.Length
$tr.Characters($pos, $d.Length).Text = $d
$pos += $d.Length
$tr.Characters($pos, $e.Length).Text = $e

# --- Paragraph 1 ---------------------------------------------------------
# "The need for Broadly " / "and Deeply Integrated " / "Middleware"
# -> "The need for broadly " / "and" / " deeply Integrated capabilities"
$r1Old = "The need for Broadly "
$r2Old = "and Deeply Integrated "
$r3Old = "Middleware"

$r1Start = $fullText.IndexOf($r1Old) + 1
$tr.Characters($r1Start, $r1Old.Length).Text = "The need for broadly "

$r2Start = $r1Start + $r1Old.Length
$tr.Characters($r2Start, $r2Old.Length).Text = "and"

$r3Start = $r2Start + "and".Length
$tr.Characters($r3Start, $r3Old.Length).Text = " deeply Integrated capabilities"
